$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching the formatting of the existing
# header cells (e.g. G1: bold, bordered, centered) by copying its format.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Add the corresponding data value in H2 (plain, unstyled numeric cell,
# matching the formatting of the other row-2 data cells like F2/G2)
$ws.Range("H2").Value = 1
